$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cells to Text format first so numeric-looking strings (e.g. "69.344.64",
# "7.10", "2.00", "0.0900") are preserved verbatim instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.344.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.679.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "683.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.679.45"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.47%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -8.71%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.301.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -10.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.680.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.330.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -9.24%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.645"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -9.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.827.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -11.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -12.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -10.00%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -8.51%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -11.91%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.59"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.12%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -10.64%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.69"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.99%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.14"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -11.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.67%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0900"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -9.73%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.941"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "166.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.69"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("B46").Value = "SuiNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -14.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "28.02"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.39%  "
